# Fix the statement_section / statement_sub_section tagging on the "cbs_6"
# balance-sheet worksheet.
#
# The Assets block (rows 2-21) was mistakenly tagged as
# "equity_liabilities" / "equity" (or "current") instead of
# "assets" / "current" | "noncurrent", while the Liabilities & Equity
# block (rows 22-47) was mistakenly tagged "assets" / "current"
# (or "noncurrent") instead of "equity_liabilities" / "current" |
# "noncurrent" | "equity".  Re-apply the correct tags to columns E and F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cbs_6")

# Assets section (rows 2-21): column E = "assets"
$ws.Range("E2:E21").Value = "assets"

# Assets - current (rows 2-11)
$ws.Range("F2:F11").Value = "current"

# Assets - non-current (rows 12-21)
$ws.Range("F12:F21").Value = "noncurrent"

# Liabilities and equity section (rows 22-47): column E = "equity_liabilities"
$ws.Range("E22:E47").Value = "equity_liabilities"

# "Total assets" row (22) sub-section
$ws.Range("F22").Value = "current"

# "Liabilities and equity" header row (23) sub-section
$ws.Range("F23").Value = "equity"

# Current liabilities (rows 24-32)
$ws.Range("F24:F32").Value = "current"

# Non-current liabilities (rows 33-40)
$ws.Range("F33:F40").Value = "noncurrent"

# Equity (rows 41-47), including "Total liabilities" / "Net assets" rows
$ws.Range("F41:F47").Value = "equity"

$wb.Save()
